$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$nl = [char]10

$ws.Range("B20").Value = "10-200 Agriculture, Industry & ConsTurnover" + $nl + "10-50 Commerce & Services"
$ws.Range("C20").Value = "< VND 20 Billionlion Agriculture, Industry & ConsTurnover" + $nl + "< VND 10 Billionlion Commerce & Services"

$ws.Range("B21").Value = "200-300 Agriculture, Industry & ConsTurnover" + $nl + "50-100 Commerce & Services"
$ws.Range("C21").Value = "VND 20 - 100 Billionlion Agriculture, Industry & ConsTurnover" + $nl + "VND 10-50 Billionlion Commerce & Services"

$ws.Range("B22").Value = ">300 Agriculture, Industry & ConsTurnover" + $nl + ">100 Commerce & Services"
$ws.Range("C22").Value = "> VND 100 Billionlion Agriculture, Industry & ConsTurnover" + $nl + "> VND 50 Billionlion Commerce & Services"
